$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Part 1: merge "...verliezer" + " opslaan." into a single run/sentence,
# then split off a new paragraph right after it that contains the new
# "In de code ..." sentence, and move the _GoBack bookmark to sit right
# after that new sentence.
# -----------------------------------------------------------------

# Step 1a: join the two runs into "...de verliezer opslaan." and insert a
# paragraph break right after the sentence (this also removes the old
# bookmark, which sat between the two runs).
$d.Content.Find.Execute(
    "telkens de verliezer opslaan.", $true, $false, $false, $false, $false,
    $true, 1, $false, "telkens de verliezer opslaan.^p", 2) | Out-Null

# Step 1b: locate the freshly created (empty) paragraph and fill it with
# the new sentence about "for each" / iterators.
$newPara = $d.Paragraphs.Item(14)
$insertPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertPoint.InsertAfter("In de code maken we gebruik van een for each lus i.p.v. een iterator uit de STL library omdat dat sneller te implementeren is en ook overzichtelijker is. Dit kan echter wel gepaard gaan met een lichte achteruitgang van de snelheid, maar dit zorgt niet voor merkbare problemen.X")

$newParaRange = $d.Paragraphs.Item(14).Range
$newParaRange.Font.Name = "Times New Roman"
$newParaRange.Font.NameAscii = "Times New Roman"
$newParaRange.Font.NameBi = "Times New Roman"
$newParaRange.Font.NameOther = "Times New Roman"
$newParaRange.Font.Size = 12

# Step 1c: re-insert the _GoBack bookmark right after the new sentence
# (just before the trailing dummy "X" marker, which keeps us away from
# the exact end-of-paragraph offset).
$p = $d.Paragraphs.Item(14)
$bmPos = $p.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Step 1d: remove the dummy "X" placeholder character.
$p2 = $d.Paragraphs.Item(14)
$delPos = $p2.Range.End - 2
$d.Range($delPos, $delPos + 1).Delete()
